$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update product name on both sheets to the new value
$newName = "4272-MS-EI-DB-DL-REC-RNI-FEE-FFConMONTHLYonLASTSUNDAY-FIFC-1-FFROP-DAILY-FIFR-1-MD-TR-1st"
$wsInput.Range("B1").Value = $newName
$wsOutput.Range("B1").Value = $newName

# Update shortname to a text value "427b" (was numeric 4272)
$wsInput.Range("B2").Value = "427b"

# Update selections: sheet1 -> B2 selected (no longer tab-selected),
# sheet2 -> B1 selected and becomes the tab-selected / active sheet.
$wsInput.Range("B2").Select() | Out-Null
$wsOutput.Activate() | Out-Null
$wsOutput.Range("B1").Select() | Out-Null
